$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I11").Value = "bb"
